$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" (column D) values look like plain numbers (e.g. "29.80",
# "0.0770") and Excel would silently coerce them to numeric values,
# dropping significant trailing zeros / formatting. Force the whole
# Price column to Text before writing so every value round-trips verbatim,
# then restore the default (Normal) style once all writes are done.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '63.361.17'
$ws.Range('E2').Value = '  +2.87%  '
$ws.Range('D3').Value = '3.487.32'
$ws.Range('E3').Value = '  +2.79%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').Value = '584.91'
$ws.Range('E5').Value = '  +1.46%  '
$ws.Range('D6').Value = '147.72'
$ws.Range('E6').Value = '  +5.15%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  +0.61%  '
$ws.Range('D9').Value = '7.71'
$ws.Range('E9').Value = '  -0.11%  '
$ws.Range('E10').Value = '  +2.94%  '
$ws.Range('E11').Value = '  +2.83%  '
$ws.Range('D12').Value = '4.083.91'
$ws.Range('E12').Value = '  +2.94%  '
$ws.Range('D13').Value = '29.80'
$ws.Range('E13').Value = '  +3.93%  '
$ws.Range('D15').Value = '3.509.12'
$ws.Range('E15').Value = '  +3.85%  '
$ws.Range('E16').Value = '  +2.17%  '
$ws.Range('D17').Value = '63.350.35'
$ws.Range('E17').Value = '  +3.25%  '
$ws.Range('D18').Value = '6.32'
$ws.Range('E18').Value = '  +3.05%  '
$ws.Range('D19').Value = '14.41'
$ws.Range('E19').Value = '  +5.58%  '
$ws.Range('D20').Value = '9.37'
$ws.Range('E20').Value = '  +4.43%  '
$ws.Range('D21').Value = '391.88'
$ws.Range('E21').Value = '  +0.52%  '
$ws.Range('D22').Value = '0.565'
$ws.Range('E22').Value = '  +1.68%  '
$ws.Range('D23').Value = '75.22'
$ws.Range('E23').Value = '  -0.30%  '
$ws.Range('E24').Value = '  -0.17%  '
$ws.Range('E25').Value = '  +5.57%  '
$ws.Range('D26').Value = '3.628.14'
$ws.Range('E26').Value = '  +3.33%  '
$ws.Range('D27').Value = '0.182'
$ws.Range('E27').Value = '  -4.20%  '
$ws.Range('D28').Value = '7.84'
$ws.Range('E28').Value = '  +8.13%  '
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('D30').Value = '8.30'
$ws.Range('E30').Value = '  +3.37%  '
$ws.Range('D31').Value = '1.47'
$ws.Range('E31').Value = '  +7.62%  '
$ws.Range('E32').Value = '  +0.81%  '
$ws.Range('D34').Value = '23.86'
$ws.Range('E34').Value = '  +2.18%  '
$ws.Range('D35').Value = '32.69'
$ws.Range('E35').Value = '  +25.86%  '
$ws.Range('E36').Value = '  +6.65%  '
$ws.Range('D37').Value = '7.14'
$ws.Range('E37').Value = '  +3.16%  '
$ws.Range('D38').Value = '171.27'
$ws.Range('E38').Value = '  +2.31%  '
$ws.Range('D39').Value = '1.58'
$ws.Range('E39').Value = '  +7.36%  '
$ws.Range('D40').Value = '3.522.29'
$ws.Range('E40').Value = '  +2.80%  '
$ws.Range('D41').Value = '0.0770'
$ws.Range('E41').Value = '  +0.79%  '
$ws.Range('D42').Value = '0.807'
$ws.Range('E42').Value = '  +3.65%  '
$ws.Range('D43').Value = '4.51'
$ws.Range('E43').Value = '  +2.30%  '
$ws.Range('D44').Value = '42.46'
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('E45').Value = '  +4.17%  '
$ws.Range('E46').Value = '  +7.51%  '
$ws.Range('D47').Value = '2.620.56'
$ws.Range('E47').Value = '  +6.24%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = '23.67'
$ws.Range('E48').Value = '  +5.87%  '
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').Value = '2.30'
$ws.Range('E49').Value = '  +13.14%  '
$ws.Range('D50').Value = '6.77'
$ws.Range('E50').Value = '  +1.17%  '
$ws.Range('E51').Value = '  +4.05%  '

$ws.Range("D2:D51").Style = "Normal"
